$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting onto the rows that will receive data / stay blank ---
# The existing blank row 19 has the plain "white, no border on date" look
# used by the freshly-added blank rows 21, 23 and 25. Grab this pattern
# before row 19 itself gets overwritten below.
$ws.Range("B19:E19").Copy()
$ws.Range("B21:E21").PasteSpecial(-4122)
$ws.Range("B23:E23").PasteSpecial(-4122)
$ws.Range("B25:E25").PasteSpecial(-4122)

# Row 16 has the "grey / date-bordered" look used by rows 18, 20, 22 and 24.
$ws.Range("B16:E16").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)
$ws.Range("B20:E20").PasteSpecial(-4122)
$ws.Range("B22:E22").PasteSpecial(-4122)
$ws.Range("B24:E24").PasteSpecial(-4122)

# Row 17 has the "white / date-bordered" look used by row 19.
$ws.Range("B17:E17").Copy()
$ws.Range("B19:E19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Fill in the values for rows 18-20 ---
$ws.Range("B18").Value = "Mercredi"
$ws.Range("C18").Value = 44579
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = "Dimentionnement et choix de composants et design du schéma électriques"

$ws.Range("B19").Value = "Dimanche"
$ws.Range("C19").Value = 44583
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = "Rédaction du rapport et design du schéma électrique"

$ws.Range("B20").Value = "Mardi"
$ws.Range("C20").Value = 44585
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = "Rédaction du rapport et corrections du schéma électrique"

# Rows 21-25 stay blank (just formatted above).

# Match the author's resulting active-cell selection.
$ws.Range("B21").Select() | Out-Null
